$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = 0.1405756097560976
$ws.Range("Z2").Value = -0.1046950969272345
$ws.Range("AB2").Value = -422.4174382225389
$ws.Range("AD2").Value = -422.4174382225389

$ws.Range("T3").Value = 0.144009756097561
$ws.Range("Z3").Value = -0.1446956540206924
$ws.Range("AB3").Value = -953.2200457609725
$ws.Range("AD3").Value = -953.2200457609725

$ws.Range("T4").Value = 0.1462634146341464
$ws.Range("Z4").Value = -0.1096439352270558
$ws.Range("AB4").Value = -702.3876576639086
$ws.Range("AD4").Value = -702.3876576639086

$ws.Range("T5").Value = 0.1500390243902439
$ws.Range("Z5").Value = -0.1146140770504377
$ws.Range("AB5").Value = -715.7955692664369
$ws.Range("AD5").Value = -715.7955692664369

$ws.Range("T6").Value = 0.1463707317073171
$ws.Range("Z6").Value = -0.1045252230381822
$ws.Range("AB6").Value = -652.9671177939568
$ws.Range("AD6").Value = -652.9671177939568

$ws.Range("T7").Value = 0.1465658536585366
$ws.Range("Z7").Value = -0.1506547063620707
$ws.Range("AB7").Value = -751.7544356918407
$ws.Range("AD7").Value = -751.7544356918407

$ws.Range("T8").Value = 0.1449658536585366
$ws.Range("Z8").Value = -0.147545094730805
$ws.Range("AB8").Value = -595.1477870917063
$ws.Range("AD8").Value = -595.1477870917063

$ws.Range("T9").Value = 0.1544
$ws.Range("Z9").Value = 0.0003559990850761401

$ws.Range("T10").Value = 0.1405756097560976
$ws.Range("Z10").Value = 0.1284680504229015
$ws.Range("AB10").Value = 518.3351116318549
$ws.Range("AD10").Value = 518.3351116318549

$ws.Range("T11").Value = 0.144009756097561
$ws.Range("Z11").Value = 0.1314260133893951
$ws.Range("AB11").Value = 865.8028559676426
$ws.Range("AD11").Value = 865.8028559676426

$ws.Range("T12").Value = 0.1462634146341464
$ws.Range("Z12").Value = 0.05484914611156014
$ws.Range("AB12").Value = 351.3679364242451
$ws.Range("AD12").Value = 351.3679364242451

$ws.Range("T13").Value = 0.1500390243902439
$ws.Range("Z13").Value = 0.1104277179693405
$ws.Range("AB13").Value = 689.6506369970006
$ws.Range("AD13").Value = 689.6506369970006

$ws.Range("T14").Value = 0.1463707317073171
$ws.Range("Z14").Value = 0.03554035381188628
$ws.Range("AB14").Value = 222.0199270509822
$ws.Range("AD14").Value = 222.0199270509822

$ws.Range("T15").Value = 0.1465658536585366
$ws.Range("Z15").Value = 0.1422008807786601
$ws.Range("AB15").Value = 709.5705502072377
$ws.Range("AD15").Value = 709.5705502072377

$ws.Range("T16").Value = 0.1449658536585366
$ws.Range("Z16").Value = 0.1814911784620183
$ws.Range("AB16").Value = 732.0749865348427
$ws.Range("AD16").Value = 732.0749865348427

$ws.Range("T17").Value = 0.1544
$ws.Range("Z17").Value = -0.002983416900589678
